$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking "Price" column (D) updates — force Text format so Excel
# doesn't silently coerce these into numbers (source data is stored as text).
$priceUpdates = [ordered]@{
    'D2' = '242.25'
    'D3' = '21.91'
    'D4' = '5.400'
    'D5' = '0.05733'
    'D6' = '3.441'
    'D7' = '6.328'
    'D8' = '0.8088'
    'D9' = '0.8565'
    'D11' = '0.07296'
    'D12' = '0.03077'
    'D13' = '0.03129'
    'D14' = '0.09380'
    'D15' = '3.949'
    'D17' = '0.04830'
    'D18' = '0.006399'
    'D19' = '0.004094'
    'D20' = '0.001000'
    'D21' = '0.0001506'
    'D22' = '3.725'
    'D23' = '2.180'
    'D24' = '0.01084'
    'D25' = '0.3264'
    'D27' = '0.0004016'
    'D40' = '0.03850'
    'D41' = '0.006678'
    'D42' = '0.1062'
    'D43' = '0.002785'
    'D44' = '0.007354'
    'D45' = '0.00005606'
    'D46' = '0.00000000753'
    'D47' = '0.5823'
    'D48' = '0.1437'
    'D49' = '0.00002108'
    'D50' = '0.01014'
}

foreach ($cellRef in $priceUpdates.Keys) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $priceUpdates[$cellRef]
}

# Plain text columns (Coin name / Link / Volume label) — assigned as-is,
# these aren't numeric-looking so no text coercion is needed.
$textUpdates = [ordered]@{
    'E13' = '12BitrueCoinBTR'
    'B18' = 'TigerCash'
    'C18' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'E18' = '17TigerCashTCH'
    'B19' = 'HotbitToken'
    'C19' = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
    'E19' = '18HotbitTokenHTB'
    'B20' = 'BitKan'
    'C20' = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
    'E20' = '19BitKanKAN'
    'B21' = 'NitroEx'
    'C21' = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
    'E21' = '20NitroExNTX'
    'B22' = 'LEO'
    'C22' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'E22' = '21LEOLEO'
    'B23' = 'BTSEToken'
    'C23' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'E23' = '22BTSETokenBTSE'
    'B24' = 'One'
    'C24' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'E24' = '23OneONEBestin24h'
}

foreach ($cellRef in $textUpdates.Keys) {
    $ws.Range($cellRef).Value = $textUpdates[$cellRef]
}
